$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 32.73910791935287
$ws.Range("C2").Value = 21.13787301585845
$ws.Range("D2").Value = 11.90132739170498
$ws.Range("E2").Value = 10.33933697746992
$ws.Range("G2").Value = 3.88736085191548
$ws.Range("I2").Value = 59.88686617569997
$ws.Range("J2").Value = 6.729235125382684
$ws.Range("L2").Value = 16.54618218790092

$ws.Range("B3").Value = 32.776308991348
$ws.Range("C3").Value = 20.83973186987473
$ws.Range("D3").Value = 11.92199023488689
$ws.Range("E3").Value = 10.35763698239853
$ws.Range("G3").Value = 3.893693232918658
$ws.Range("I3").Value = 58.84916305600157
$ws.Range("J3").Value = 6.712316352058839
$ws.Range("L3").Value = 16.55879517701753

$ws.Range("B4").Value = 32.81275502645713
$ws.Range("C4").Value = 20.66189095383691
$ws.Range("D4").Value = 11.93670869575271
$ws.Range("E4").Value = 10.36951598078669
$ws.Range("G4").Value = 3.897769505781398
$ws.Range("I4").Value = 58.20450931094133
$ws.Range("J4").Value = 6.701702831340215
$ws.Range("L4").Value = 16.5698669429226

$ws.Range("B5").Value = 32.83100876448321
$ws.Range("C5").Value = 20.59081437087639
$ws.Range("D5").Value = 11.94321644243216
$ws.Range("E5").Value = 10.37451899172855
$ws.Range("G5").Value = 3.899478191902513
$ws.Range("I5").Value = 57.94011763452054
$ws.Range("J5").Value = 6.697320752477663
$ws.Range("L5").Value = 16.57521326215322

$ws.Range("B6").Value = 32.83424462754301
$ws.Range("C6").Value = 20.57909876522485
$ws.Range("D6").Value = 11.94432781394608
$ws.Range("E6").Value = 10.37535955351036
$ws.Range("G6").Value = 3.89976479864112
$ws.Range("I6").Value = 57.89611931928974
$ws.Range("J6").Value = 6.696589660519334
$ws.Range("L6").Value = 16.57615134632719

$ws.Range("B7").Value = 32.81298745658004
$ws.Range("C7").Value = 20.66092663492688
$ws.Range("D7").Value = 11.93679439839871
$ws.Range("E7").Value = 10.36958279549955
$ws.Range("G7").Value = 3.897792356762376
$ws.Range("I7").Value = 58.20095021380786
$ws.Range("J7").Value = 6.701643964244873
$ws.Range("L7").Value = 16.56993566951174

$ws.Range("B8").Value = 32.74910038995573
$ws.Range("C8").Value = 21.03404132427228
$ws.Range("D8").Value = 11.9080297613208
$ws.Range("E8").Value = 10.34551378691908
$ws.Range("G8").Value = 3.889505358075611
$ws.Range("I8").Value = 59.53073891889704
$ws.Range("J8").Value = 6.723448101085851
$ws.Range("L8").Value = 16.54983929631168

$ws.Range("B9").Value = 32.73251113278619
$ws.Range("C9").Value = 21.80335794952994
$ws.Range("D9").Value = 11.86778283379565
$ws.Range("E9").Value = 10.30338684440447
$ws.Range("G9").Value = 3.874735855847342
$ws.Range("I9").Value = 62.07059486452516
$ws.Range("J9").Value = 6.764433948211917
$ws.Range("L9").Value = 16.53692256434059

$ws.Range("B10").Value = 32.78743870146994
$ws.Range("C10").Value = 22.38635004093894
$ws.Range("D10").Value = 11.84812242383956
$ws.Range("E10").Value = 10.27549034648489
$ws.Range("G10").Value = 3.864771232823129
$ws.Range("I10").Value = 63.88452460431567
$ws.Range("J10").Value = 6.793489856169664
$ws.Range("L10").Value = 16.54369732559697

$ws.Range("B11").Value = 32.82712586863737
$ws.Range("C11").Value = 22.6543178279533
$ws.Range("D11").Value = 11.84134143545017
$ws.Range("E11").Value = 10.26345457184539
$ws.Range("G11").Value = 3.860426995267983
$ws.Range("I11").Value = 64.69627958401756
$ws.Range("J11").Value = 6.806484013578262
$ws.Range("L11").Value = 16.55033062920607

$ws.Range("B12").Value = 32.84427389187501
$ws.Range("C12").Value = 22.75609542493532
$ws.Range("D12").Value = 11.8390854529178
$ws.Range("E12").Value = 10.25899041948281
$ws.Range("G12").Value = 3.858808800751699
$ws.Range("I12").Value = 65.00157292650506
$ws.Range("J12").Value = 6.811372800808183
$ws.Range("L12").Value = 16.55335425823908

$ws.Range("B13").Value = 32.84048644334106
$ws.Range("C13").Value = 22.73416374112422
$ws.Range("D13").Value = 11.83955743747317
$ws.Range("E13").Value = 10.25994770381495
$ws.Range("G13").Value = 3.859156116882651
$ws.Range("I13").Value = 64.93591861641737
$ws.Range("J13").Value = 6.810321322936497
$ws.Range("L13").Value = 16.55268029194684

$ws.Range("B14").Value = 32.82849415095021
$ws.Range("C14").Value = 22.66268556634896
$ws.Range("D14").Value = 11.84114958020333
$ws.Range("E14").Value = 10.26308543209058
$ws.Range("G14").Value = 3.860293328354541
$ws.Range("I14").Value = 64.72143930577552
$ws.Range("J14").Value = 6.806886852859131
$ws.Range("L14").Value = 16.55056912388014

$ws.Range("B15").Value = 32.82142464023489
$ws.Range("C15").Value = 22.61893997769675
$ws.Range("D15").Value = 11.84216544734599
$ws.Range("E15").Value = 10.26501954501519
$ws.Range("G15").Value = 3.860993395141843
$ws.Range("I15").Value = 64.58978600989256
$ws.Range("J15").Value = 6.804779006784674
$ws.Range("L15").Value = 16.54934264281886

$ws.Range("B16").Value = 32.7851413938881
$ws.Range("C16").Value = 22.3688852912555
$ws.Range("D16").Value = 11.84860916869045
$ws.Range("E16").Value = 10.27629003294489
$ws.Range("G16").Value = 3.865058913984839
$ws.Range("I16").Value = 63.83118972630598
$ws.Range("J16").Value = 6.792636187213681
$ws.Range("L16").Value = 16.54333538330761

$ws.Range("B17").Value = 32.7666527208408
$ws.Range("C17").Value = 22.21612543771091
$ws.Range("D17").Value = 11.8531167046215
$ws.Range("E17").Value = 10.28337132796867
$ws.Range("G17").Value = 3.867601126657606
$ws.Range("I17").Value = 63.3622553484145
$ws.Range("J17").Value = 6.785129940971912
$ws.Range("L17").Value = 16.54056056254623

$ws.Range("B18").Value = 32.75740204970498
$ws.Range("C18").Value = 22.12852745364939
$ws.Range("D18").Value = 11.85591282999476
$ws.Range("E18").Value = 10.28750593908035
$ws.Range("G18").Value = 3.869081119194723
$ws.Range("I18").Value = 63.09128416810283
$ws.Range("J18").Value = 6.780791488254359
$ws.Range("L18").Value = 16.53929879088013

$ws.Range("B19").Value = 32.75450733842246
$ws.Range("C19").Value = 22.09891673541552
$ws.Range("D19").Value = 11.85689447630519
$ws.Range("E19").Value = 10.28891645133877
$ws.Range("G19").Value = 3.869585280804816
$ws.Range("I19").Value = 62.9993282508129
$ws.Range("J19").Value = 6.779318936109588
$ws.Range("L19").Value = 16.53892894490227

$ws.Range("B20").Value = 32.76847763330395
$ws.Range("C20").Value = 22.23236019245271
$ws.Range("D20").Value = 11.85261580091817
$ws.Range("E20").Value = 10.28261113698175
$ws.Range("G20").Value = 3.867328665619763
$ws.Range("I20").Value = 63.41230512101438
$ws.Range("J20").Value = 6.785931170828609
$ws.Range("L20").Value = 16.54082134678245

$ws.Range("B21").Value = 32.83195903017305
$ws.Range("C21").Value = 22.68367292152777
$ws.Range("D21").Value = 11.84067345981995
$ws.Range("E21").Value = 10.26216127135243
$ws.Range("G21").Value = 3.859958574537006
$ws.Range("I21").Value = 64.78449541744499
$ws.Range("J21").Value = 6.807896500910769
$ws.Range("L21").Value = 16.55117532957909

$ws.Range("B22").Value = 32.88580148072856
$ws.Range("C22").Value = 22.98036209218345
$ws.Range("D22").Value = 11.83468639099196
$ws.Range("E22").Value = 10.24934102363837
$ws.Range("G22").Value = 3.855298312247501
$ws.Range("I22").Value = 65.66898342044557
$ws.Range("J22").Value = 6.822066806514387
$ws.Range("L22").Value = 16.5609253711339

$ws.Range("B23").Value = 32.85593333857567
$ws.Range("C23").Value = 22.82188545896586
$ws.Range("D23").Value = 11.83771519492341
$ws.Range("E23").Value = 10.25613376176891
$ws.Range("G23").Value = 3.857771348110694
$ws.Range("I23").Value = 65.19809715510134
$ws.Range("J23").Value = 6.81452069532397
$ws.Range("L23").Value = 16.55544835771899

$ws.Range("B24").Value = 32.76764829644426
$ws.Range("C24").Value = 22.22501974977407
$ws.Range("D24").Value = 11.85284162200574
$ws.Range("E24").Value = 10.28295462155626
$ws.Range("G24").Value = 3.867451787796417
$ws.Range("I24").Value = 63.3896818871079
$ws.Range("J24").Value = 6.785569006511291
$ws.Range("L24").Value = 16.54070240750357

$ws.Range("B25").Value = 32.72526383418491
$ws.Range("C25").Value = 21.59173788602339
$ws.Range("D25").Value = 11.87693450297723
$ws.Range("E25").Value = 10.31424421351259
$ws.Range("G25").Value = 3.878574516541977
$ws.Range("I25").Value = 61.39189471541899
$ws.Range("J25").Value = 6.753533432594729
$ws.Range("L25").Value = 16.5375679012399
